$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and this card has been revealed.
</p><p><b>Miracle:</b> At the start of the Investigator Phase, if you have a Light Source, a Bladed Weapon, and an Evidence Unique Item, and you are sharing a space with another investigator, you may reveal this card. That investigator may discard a Wounded or Insane condition. This miracle may only happen once per game.</p>"
$ws.Range("D3").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and this card has been revealed.</p>
<p><b>Execution Attempt:</b> At the start of the Investigator Phase, if you have a Bladed Weapon and you are in a space with exactly one other investigator, you may reveal this card. If that investigator is not Wounded, he becomes Wounded and discards all facedown Damage. If he is Wounded, that investigator is eliminated. The game does not end as a result of that investigator’s elimination. This execution attempt may only happen once per game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>"
$ws.Range("E3").Value = ""
$ws.Range("D4").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have 6 or more items.</p><p><em>Item:</em> The word `"Items`" is a collective term that refers to both Common Items and Unique Items. Spells are not items. (Spells and Items are possessions)</p>"
$ws.Range("E4").Value = "Item: The word “Items” is a collective term that refers to both Common Items and Unique Items. Spells are not items. (Spells and Items are possessions)"
$ws.Range("D5").Value = "<p>You do not win the game as normal. Instead, you only win if the investigation is complete and this card has been revealed.
</p>
<p><b>Conflagration:</b> At the start of the Investigator Phase, if 6 or more rooms contain Fire, reveal this card. While this card is revealed, at the beginning of each Mythos phase, fire spreads twice.</p>
<p><b>Set Fire (Action):</b> If you have a Light Source, place a Fire in your space or an adjacent space. Select the “Set Fire” option in the app.</p>
<p><b>Extinguish Fire (Action):</b> Test Agility. For each success, discard Fire from your space or a space you move into later during the round. Doesn’t do Damage.</p>
<p>If group disagrees about how a fire should spread, a random investigator decides. Roll a die, draw cards, use an app, or do something else.</p>"
$ws.Range("E5").Value = "• Set Fire (Action): If you have a Light Source, place a Fire in your space or an adjacent space. Select the “Set Fire” option in the app.
• Extinguish Fire (Action): Test Agility. For each success, discard Fire from your space or a space you move into later during the round. Doesn’t do Damage.
• If group disagrees about how a fire should spread, a random investigator decides. Roll a die, draw cards, use an app, or do something else.
"
$ws.Range("D6").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and this card has been revealed.
</p>
<p><b>Martyrdom:</b> After you have been eliminated, reveal this card. The game does not end as a result of your elimination.</p><p>When eliminated, you drop all possessions onto your space and remove your figure from the board. You cannot take actions.
</p><p>Normally, eliminations cause the eliminated player to lose and the group to lose after the next Investigation phase. This card overrides that rule.</p>"
$ws.Range("E6").Value = "•	When eliminated, you drop all possessions onto your space and remove your figure from the board. You cannot take actions.
•	Normally, eliminations cause the eliminated player to lose and the group to lose after the next Investigation phase. This card overrides that rule."
$ws.Range("D7").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and this card has been revealed.
</p>
<p><b>Bloodletting:</b> At the start of the Investigator Phase, if you have a Bladed Weapon and there is at least one other investigator on your space, you may reveal this card and choose one investigator on your space. That investigator suffers one facedown Damage, and you discard one Horror. Other investigators within range flip one Horror.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>
<p><em>`"Within Range:`"</em> Means up to 3 spaces away. Cannot count through walls, doors, or impassable borders, unless an effect says otherwise.</p>"
$ws.Range("E7").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession.
•	“Within Range:” Means up to 3 spaces away. Cannot count through walls, doors, or impassable borders, unless an effect says otherwise."
$ws.Range("D8").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and one or more other investigators are Insane. Otherwise, you lose the game.</p><p><b>Push (Action):</b> Choose target monster/investigator in your space, choose destination, and determine willingness. Skip test if willing. Determine test difficulty: monsters use Brawn; investigators roll Strength and add 1 to test result. Roll Strength test. If you succeed, push target to chosen space. You may also move to that space.</p>"
$ws.Range("E8").Value = "•	Push (Action): Choose target monster/investigator in your space, choose destination, and determine willingness. Skip test if willing. Determine test difficulty: monsters use Brawn; investigators roll Strength and add 1 to test result. Roll Strength test. If you succeed, push target to chosen space. You may also move to that space."
$ws.Range("D9").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and one or more other investigators are Wounded. Otherwise, you lose the game.</p>"
$ws.Range("D10").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and there are no Search tokens on the board. Otherwise, you lose the game.</p><p>Search tokens are question marks:    ?</p>
<p>You do not need to clear the Interact, Explore, Sight, and Person tokens.</p>"
$ws.Range("E10").Value = "•	Search tokens are question marks:    ?
•	You do not need to clear the Interact, Explore, Sight, and Person tokens."
$ws.Range("D11").Value = "<p>You cannot perform any single action more than once each round. You win or lose the game as normal.</p><p>You may still perform 2 actions each round, but they must each be different actions.</p>"
$ws.Range("E11").Value = "You may still perform 2 actions each round, but they must each be different actions."
$ws.Range("D12").Value = "<p><em>No effect.</em> You win or lose the game as normal.</p><p>Remember, you cannot reveal the back of your Insanity card.</p>"
$ws.Range("E12").Value = "Remember, you cannot reveal the back of your Insanity. "
$ws.Range("D13").Value = "<p>You do not win the game as normal. Instead, you win if the investigation is complete and this card has been revealed.
</p>
<p><b>Burn the Evidence:</b> At the start of the Investigator Phase, if you possess two Evidence Unique Items while standing on or adjacent to a space that contains Fire, you may reveal this card. Discard two Evidence Unique items from your inventory.</p><p><b>Set Fire (Action):</b> If you have a Light Source, place a Fire in your space or an adjacent space. Select the `"Set Fire`" option in the app.</p>
<p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>"
$ws.Range("E13").Value = "•	Set Fire (Action): If you have a Light Source, place a Fire in your space or an adjacent space. Select the “Set Fire” option in the app.
•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession."
$ws.Range("D14").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have 1 or more Spells. Otherwise, you lose the game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession. (Spells are possessions)</p>"
$ws.Range("E14").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession. (Spells are possessions)"
$ws.Range("D15").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have more Spells than the investigator to your right. Otherwise, you lose the game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession. (Spells are possessions)</p>"
$ws.Range("E15").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession. (Spells are possessions)"
$ws.Range("D16").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have more Evidence than the investigator to your right. Otherwise, you lose the game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>"
$ws.Range("E16").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession."
$ws.Range("D17").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have more Items than the investigator to your left. Otherwise, you lose the game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>
<p><em>Item:</em> The word `"Items`" is a collective term that refers to both Common Items and Unique Items. Spells are not Items.</p>"
$ws.Range("E17").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession.
•	Item: The word “Items” is a collective term that refers to both Common Items and Unique Items. Spells are not Items."
$ws.Range("D18").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and you have 1 or more Bladed Weapons or Light Sources.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession.</p>"
$ws.Range("E18").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession."
$ws.Range("D19").Value = "<p>You cannot use Actions defined on your character sheet. You win or lose the game as normal.</p><p>Other non-Action character abilities take place normally (such as gaining extra clues when searching, etc).</p><p>Remember, you cannot reveal the back of your insanity.</p>"
$ws.Range("E19").Value = "•	Other non-Action character abilities take place normally (such as gaining extra clues when searching, etc).
•	Remember, you cannot reveal the back of your insanity."
$ws.Range("D20").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and this card has been revealed.
</p>
<p><b>Horrific Offering:</b> At the start of the Investigator Phase, if you have a Bladed Weapon and 1 or more Spells, and there is at least one other investigator on your space, you may reveal this card and choose one investigator on your space. That investigator suffers one Damage and becomes Insane. (If the investigator was already Insane, the investigator is eliminated as normal.) This offering may only happen once per game.</p><p><b>Steal (Action):</b> Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent's, you take one possession. (Spells are possessions)</p>"
$ws.Range("E20").Value = "•	Steal (Action): Pick an investigator in your space. Pick Strength, Agility, or Observation. Each player rolls that test. For each net success over your opponent’s, you take one possession. (Spells are possessions)"
